$changes = @(
    @{ Row = 2; Col = 4; Value = "30.352.61"; ForceText = $false },
    @{ Row = 2; Col = 5; Value = "  -2.67%  "; ForceText = $false },
    @{ Row = 3; Col = 4; Value = "1.939.00"; ForceText = $false },
    @{ Row = 3; Col = 5; Value = "  -2.68%  "; ForceText = $false },
    @{ Row = 4; Col = 4; Value = "1.002"; ForceText = $true },
    @{ Row = 4; Col = 5; Value = "  +0.12%  "; ForceText = $false },
    @{ Row = 5; Col = 4; Value = "250.36"; ForceText = $true },
    @{ Row = 5; Col = 5; Value = "  -1.61%  "; ForceText = $false },
    @{ Row = 6; Col = 4; Value = "0.7260"; ForceText = $true },
    @{ Row = 6; Col = 5; Value = "  -6.64%  "; ForceText = $false },
    @{ Row = 7; Col = 5; Value = "  +0.02%  "; ForceText = $false },
    @{ Row = 8; Col = 4; Value = "0.3332"; ForceText = $true },
    @{ Row = 8; Col = 5; Value = "  -3.96%  "; ForceText = $false },
    @{ Row = 9; Col = 4; Value = "28.44"; ForceText = $true },
    @{ Row = 9; Col = 5; Value = "  +2.68%  "; ForceText = $false },
    @{ Row = 10; Col = 4; Value = "0.07262"; ForceText = $true },
    @{ Row = 10; Col = 5; Value = "  +3.28%  "; ForceText = $false },
    @{ Row = 11; Col = 4; Value = "0.8109"; ForceText = $true },
    @{ Row = 11; Col = 5; Value = "  -4.02%  "; ForceText = $false },
    @{ Row = 12; Col = 4; Value = "0.08096"; ForceText = $true },
    @{ Row = 12; Col = 5; Value = "  -0.96%  "; ForceText = $false },
    @{ Row = 13; Col = 4; Value = "1.937.35"; ForceText = $false },
    @{ Row = 13; Col = 5; Value = "  -2.69%  "; ForceText = $false },
    @{ Row = 14; Col = 4; Value = "5.476"; ForceText = $true },
    @{ Row = 14; Col = 5; Value = "  -2.66%  "; ForceText = $false },
    @{ Row = 15; Col = 4; Value = "94.62"; ForceText = $true },
    @{ Row = 15; Col = 5; Value = "  -5.96%  "; ForceText = $false },
    @{ Row = 16; Col = 5; Value = "  -1.95%  "; ForceText = $false },
    @{ Row = 17; Col = 4; Value = "30.357.92"; ForceText = $false },
    @{ Row = 17; Col = 5; Value = "  -2.70%  "; ForceText = $false },
    @{ Row = 18; Col = 4; Value = "0.000008237"; ForceText = $true },
    @{ Row = 18; Col = 5; Value = "  +2.69%  "; ForceText = $false },
    @{ Row = 19; Col = 4; Value = "250.89"; ForceText = $true },
    @{ Row = 19; Col = 5; Value = "  -7.90%  "; ForceText = $false },
    @{ Row = 20; Col = 4; Value = "5.918"; ForceText = $true },
    @{ Row = 20; Col = 5; Value = "  +0.71%  "; ForceText = $false },
    @{ Row = 21; Col = 4; Value = "2.192.16"; ForceText = $false },
    @{ Row = 21; Col = 5; Value = "  -2.69%  "; ForceText = $false },
    @{ Row = 22; Col = 5; Value = "  -0.06%  "; ForceText = $false },
    @{ Row = 23; Col = 4; Value = "1.003"; ForceText = $true },
    @{ Row = 23; Col = 5; Value = "  +0.29%  "; ForceText = $false },
    @{ Row = 24; Col = 4; Value = "6.948"; ForceText = $true },
    @{ Row = 24; Col = 5; Value = "  -1.50%  "; ForceText = $false },
    @{ Row = 25; Col = 4; Value = "9.771"; ForceText = $true },
    @{ Row = 25; Col = 5; Value = "  -1.99%  "; ForceText = $false },
    @{ Row = 26; Col = 4; Value = "163.01"; ForceText = $true },
    @{ Row = 26; Col = 5; Value = "  -1.67%  "; ForceText = $false },
    @{ Row = 27; Col = 4; Value = "2.391"; ForceText = $true },
    @{ Row = 27; Col = 5; Value = "  +1.40%  "; ForceText = $false },
    @{ Row = 28; Col = 4; Value = "19.32"; ForceText = $true },
    @{ Row = 28; Col = 5; Value = "  -2.77%  "; ForceText = $false },
    @{ Row = 29; Col = 5; Value = "  -7.23%  "; ForceText = $false },
    @{ Row = 30; Col = 4; Value = "1.570"; ForceText = $true },
    @{ Row = 30; Col = 5; Value = "  -1.77%  "; ForceText = $false },
    @{ Row = 31; Col = 4; Value = "1.352"; ForceText = $true },
    @{ Row = 31; Col = 5; Value = "  -1.25%  "; ForceText = $false },
    @{ Row = 32; Col = 4; Value = "4.439"; ForceText = $true },
    @{ Row = 32; Col = 5; Value = "  -3.23%  "; ForceText = $false },
    @{ Row = 33; Col = 4; Value = "4.201"; ForceText = $true },
    @{ Row = 33; Col = 5; Value = "  -5.13%  "; ForceText = $false },
    @{ Row = 34; Col = 4; Value = "0.05190"; ForceText = $true },
    @{ Row = 34; Col = 5; Value = "  -0.96%  "; ForceText = $false },
    @{ Row = 35; Col = 4; Value = "1.292"; ForceText = $true },
    @{ Row = 35; Col = 5; Value = "  +6.53%  "; ForceText = $false },
    @{ Row = 36; Col = 4; Value = "0.7513"; ForceText = $true },
    @{ Row = 36; Col = 5; Value = "  -4.28%  "; ForceText = $false },
    @{ Row = 37; Col = 4; Value = "2.750"; ForceText = $true },
    @{ Row = 37; Col = 5; Value = "  -0.29%  "; ForceText = $false },
    @{ Row = 38; Col = 4; Value = "0.01978"; ForceText = $true },
    @{ Row = 38; Col = 5; Value = "  -1.20%  "; ForceText = $false },
    @{ Row = 39; Col = 4; Value = "2.837"; ForceText = $true },
    @{ Row = 39; Col = 5; Value = "  -1.95%  "; ForceText = $false },
    @{ Row = 40; Col = 4; Value = "80.96"; ForceText = $true },
    @{ Row = 40; Col = 5; Value = "  +1.36%  "; ForceText = $false },
    @{ Row = 41; Col = 4; Value = "6.526"; ForceText = $true },
    @{ Row = 41; Col = 5; Value = "  -2.41%  "; ForceText = $false },
    @{ Row = 42; Col = 4; Value = "0.4541"; ForceText = $true },
    @{ Row = 42; Col = 5; Value = "  -2.66%  "; ForceText = $false },
    @{ Row = 43; Col = 4; Value = "2.038"; ForceText = $true },
    @{ Row = 43; Col = 5; Value = "  -2.84%  "; ForceText = $false },
    @{ Row = 44; Col = 4; Value = "0.8482"; ForceText = $true },
    @{ Row = 44; Col = 5; Value = "  -0.73%  "; ForceText = $false },
    @{ Row = 45; Col = 5; Value = "  -0.06%  "; ForceText = $false },
    @{ Row = 46; Col = 4; Value = "102.27"; ForceText = $true },
    @{ Row = 46; Col = 5; Value = "  -2.14%  "; ForceText = $false },
    @{ Row = 47; Col = 4; Value = "9.827"; ForceText = $true },
    @{ Row = 47; Col = 5; Value = "  -1.39%  "; ForceText = $false },
    @{ Row = 48; Col = 4; Value = "7.454"; ForceText = $true },
    @{ Row = 48; Col = 5; Value = "  -2.53%  "; ForceText = $false },
    @{ Row = 49; Col = 4; Value = "36.84"; ForceText = $true },
    @{ Row = 49; Col = 5; Value = "  -1.14%  "; ForceText = $false },
    @{ Row = 50; Col = 4; Value = "0.4190"; ForceText = $true },
    @{ Row = 50; Col = 5; Value = "  -2.62%  "; ForceText = $false },
    @{ Row = 51; Col = 4; Value = "0.06056"; ForceText = $true },
    @{ Row = 51; Col = 5; Value = "  +1.33%  "; ForceText = $false },
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($change in $changes) {
    $cell = $ws.Cells.Item($change.Row, $change.Col)
    if ($change.ForceText) {
        # The new value looks like a plain number (e.g. "250.36"); without
        # forcing text format Excel would silently convert it to a numeric
        # cell, losing the original text representation. Format as Text,
        # assign, then restore the default "Normal" style so no stray
        # style index is left behind on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $change.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $change.Value
    }
}
